$d = $word.ActiveDocument

$replacements = @(
    @("836×2=1672", "759×9=6831"),
    @("913×9=8217", "827×9=7443"),
    @("952×4=3808", "505×9=4545"),
    @("751×4=3004", "774×5=3870"),
    @("900×3=2700", "120×8=960"),
    @("150×7=1050", "692×9=6228"),
    @("701×7=4907", "306×4=1224"),
    @("113×6=678",  "611×4=2444"),
    @("399×6=2394", "240×4=960"),
    @("588×5=2940", "870×7=6090"),
    @("160×2=320",  "323×4=1292"),
    @("673×4=2692", "596×6=3576"),
    @("424×7=2968", "656×6=3936"),
    @("239×5=1195", "529×9=4761"),
    @("386×4=1544", "947×3=2841"),
    @("910×2=1820", "845×9=7605"),
    @("580×6=3480", "331×3=993"),
    @("306×2=612",  "354×3=1062"),
    @("747×4=2988", "975×2=1950"),
    @("748×8=5984", "104×2=208"),
    @("872×7=6104", "793×8=6344"),
    @("744×9=6696", "510×5=2550"),
    @("310×3=930",  "778×5=3890"),
    @("674×9=6066", "423×7=2961"),
    @("350×8=2800", "638×6=3828")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
